$d = $word.ActiveDocument

# --- 1. Add the three new character styles -------------------------------

$gaNStyle = $d.Styles.Add("GaNStyle", 2)
$gaNStyle.Font.Name = "Calibri"
$gaNStyle.Font.Size = 14

$gaNParagraph = $d.Styles.Add("GaNParagraph", 2)
$gaNParagraph.Font.Name = "Calibri"
$gaNParagraph.Font.Size = 10

$gaNLinks = $d.Styles.Add("GaNLinks", 2)
$gaNLinks.Font.Name = "Calibri"
$gaNLinks.Font.Bold = $true
$gaNLinks.Font.Color = 8388608
$gaNLinks.Font.Size = 9.5
$gaNLinks.Font.Underline = 1

# --- 2. Apply GaNStyle to every "Datas das campanhas..." run -------------

$datesText = "Datas das campanhas de 2022 que usam Constelação de Órion: 16 a 25 de janeiro, 14 a 23 de fevereiro, 14 a 24 de março"

$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute($datesText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
while ($found) {
    $range.Style = "GaNStyle"
    $range.Collapse(0)
    $range.End = $d.Content.End
    $found = $range.Find.Execute($datesText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
}

# --- 3. Apply GaNParagraph to the "Está a participar..." paragraph run ---

$participarText = "Está a participar numa campanha global para observar e registar as estrelas mais fracas visíveis como forma de medir a poluição luminosa num determinado local. Localizando e observando a  Constelação de Órion no céu noturno e,  comparando-a com cartas estelares, pessoas de todo o mundo aprenderão  como as luzes da sua comunidade contribuem para a poluição luminosa. As suas contribuições para a base de dados on-line irão documentar a visibilidade do céu noturno em todo o mundo."

$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute($participarText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
while ($found) {
    $range.Style = "GaNParagraph"
    $range.Collapse(0)
    $range.End = $d.Content.End
    $found = $range.Find.Execute($participarText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
}

# --- 4. Apply GaNLinks to the "por Jenik Hollan..." run -------------------

$jenikText = "por Jenik Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."

$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute($jenikText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
while ($found) {
    $range.Style = "GaNLinks"
    $range.Collapse(0)
    $range.End = $d.Content.End
    $found = $range.Find.Execute($jenikText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
}

Write-Output "done"
